$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Sprint 2 tasks (rows 10-13) got marked complete after the micro-ui demo
# started working; rows 10-12 also got flagged as assigned to a sprint.
$ws.Range("H10:H13").Value = "Complete"
$ws.Range("I10:I12").Value = "Yes"

# Cursor ended up on D10 (one row up from where it was before, since the
# frozen/top-left-cell scroll reset back to the top of the sheet).
$ws.Range("D10").Select()
